# Feature: add arrows (arrow_n). Fixed bugs, removed unnecessary code.
#
# The "meta" worksheet holds a key/value table of chart options (tab,
# title, x_title, y_title, forecast_x, hline_bold, ...) in columns A/B,
# terminated by a blank placeholder row. This change introduces a new
# "style" = "default" option: a row is inserted above the old blank
# placeholder row (row 7), which pushes that blank row down to row 8,
# and the new row 7 is filled in with the style option.
$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("meta")

$ws1.Rows.Item(7).Insert()
$ws1.Range("A7").Value = "style"
$ws1.Range("B7").Value = "default"

# The data worksheet got a tiny floating-point precision correction on
# one of its forecast values (row 5, column C).
$ws2 = $wb.Worksheets.Item("forecasted-dashed-lines")
$ws2.Range("C5").Value = 95.659395447989
